{"js": "// Replace the 100 arithmetic answers in the single table, in document\n// (row-major) order, with their updated values. Each table cell holds\n// exactly one paragraph with one run/one text value, so we can address\n// cells positionally via Table.getCell(row, col).\nconst newValues = [\"25+62=87\", \"51-18=33\", \"58+0=58\", \"66-65=1\", \"60+7=67\", \"57+5=62\", \"55-1=54\", \"37+10=47\", \"25+59=84\", \"69-61=8\", \"83+4=87\", \"3-1=2\", \"49+20=69\", \"14+8=22\", \"36+53=89\", \"26+45=71\", \"77-76=1\", \"66+7=73\", \"39+9=48\", \"84-53=31\", \"64-52=12\", \"49+48=97\", \"36+9=45\", \"86-55=31\", \"48+11=59\", \"9+72=81\", \"85-73=12\", \"78+6=84\", \"56-6=50\", \"53+27=80\", \"28+12=40\", \"73-14=59\", \"22-20=2\", \"6+42=48\", \"60-17=43\", \"98-97=1\", \"35+60=95\", \"72+4=76\", \"28-16=12\", \"55+20=75\", \"26+57=83\", \"41-21=20\", \"80-26=54\", \"20+46=66\", \"81-77=4\", \"49-34=15\", \"88-28=60\", \"7+46=53\", \"68+6=74\", \"7+46=53\", \"20+46=66\", \"47-31=16\", \"9-3=6\", \"58-45=13\", \"49+23=72\", \"56-31=25\", \"55+2=57\", \"92-16=76\", \"83+0=83\", \"28+38=66\", \"47-30=17\", \"79+19=98\", \"3+74=77\", \"50-18=32\", \"87-29=58\", \"43+15=58\", \"95-18=77\", \"44-13=31\", \"8+78=86\", \"70+12=82\", \"93-2=91\", \"18+43=61\", \"4+22=26\", \"28+49=77\", \"94-9=85\", \"79-27=52\", \"17+49=66\", \"66+5=71\", \"71-59=12\", \"1+96=97\", \"78+8=86\", \"95-89=6\", \"22-12=10\", \"48-16=32\", \"92-10=82\", \"78-51=27\", \"23+42=65\", \"65-12=53\", \"54-9=45\", \"4+68=72\", \"2+18=20\", \"96-69=27\", \"26+1=27\", \"13+51=64\", \"59+37=96\", \"88-12=76\", \"94-68=26\", \"12-12=0\", \"74-30=44\", \"3+23=26\"];\n\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nlet idx = 0;\nfor (let r = 0; r < rows.items.length; r++) {\n  const row = rows.items[r];\n  row.cells.load(\"items\");\n  await context.sync();\n  const cells = row.cells.items;\n  for (let c = 0; c < cells.length; c++) {\n    if (idx >= newValues.length) break;\n    cells[c].value = newValues[idx];\n    idx++;\n  }\n}\nawait context.sync();\n", "ps1": "# Replace the 100 arithmetic answers in the single table, in document\n# (row-major) order, with their updated values. Each table cell holds\n# exactly one paragraph with one run/one text value, so cells are\n# addressed positionally via Table.Cell(row, col).\n$newValues = @(\"25+62=87\",\"51-18=33\",\"58+0=58\",\"66-65=1\",\"60+7=67\",\"57+5=62\",\"55-1=54\",\"37+10=47\",\"25+59=84\",\"69-61=8\",\"83+4=87\",\"3-1=2\",\"49+20=69\",\"14+8=22\",\"36+53=89\",\"26+45=71\",\"77-76=1\",\"66+7=73\",\"39+9=48\",\"84-53=31\",\"64-52=12\",\"49+48=97\",\"36+9=45\",\"86-55=31\",\"48+11=59\",\"9+72=81\",\"85-73=12\",\"78+6=84\",\"56-6=50\",\"53+27=80\",\"28+12=40\",\"73-14=59\",\"22-20=2\",\"6+42=48\",\"60-17=43\",\"98-97=1\",\"35+60=95\",\"72+4=76\",\"28-16=12\",\"55+20=75\",\"26+57=83\",\"41-21=20\",\"80-26=54\",\"20+46=66\",\"81-77=4\",\"49-34=15\",\"88-28=60\",\"7+46=53\",\"68+6=74\",\"7+46=53\",\"20+46=66\",\"47-31=16\",\"9-3=6\",\"58-45=13\",\"49+23=72\",\"56-31=25\",\"55+2=57\",\"92-16=76\",\"83+0=83\",\"28+38=66\",\"47-30=17\",\"79+19=98\",\"3+74=77\",\"50-18=32\",\"87-29=58\",\"43+15=58\",\"95-18=77\",\"44-13=31\",\"8+78=86\",\"70+12=82\",\"93-2=91\",\"18+43=61\",\"4+22=26\",\"28+49=77\",\"94-9=85\",\"79-27=52\",\"17+49=66\",\"66+5=71\",\"71-59=12\",\"1+96=97\",\"78+8=86\",\"95-89=6\",\"22-12=10\",\"48-16=32\",\"92-10=82\",\"78-51=27\",\"23+42=65\",\"65-12=53\",\"54-9=45\",\"4+68=72\",\"2+18=20\",\"96-69=27\",\"26+1=27\",\"13+51=64\",\"59+37=96\",\"88-12=76\",\"94-68=26\",\"12-12=0\",\"74-30=44\",\"3+23=26\")\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n\n$idx = 0\nfor ($r = 1; $r -le $rows; $r++) {\n  for ($c = 1; $c -le $cols; $c++) {\n    if ($idx -ge $newValues.Length) { break }\n    $cell = $t.Cell($r, $c)\n    $cell.Range.Text = $newValues[$idx]\n    $idx++\n  }\n}\n"}
